$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.13541
$ws.Range("H2").Value = 0.40623
$ws.Range("I2").Value = 0.1064658135528677
$ws.Range("J2").Value = 0.1064658135528677
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.048574
$ws.Range("N2").Value = 0.145722
$ws.Range("O2").Value = 0.005173961045635648
$ws.Range("P2").Value = 0.005173961045635648
$ws.Range("Q2").Value = 0.00657740534
$ws.Range("R2").Value = 0.05919664806
$ws.Range("S2").Value = 0.0005508499720144453
$ws.Range("T2").Value = 0.0005508499720144453

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.13541
$ws.Range("H3").Value = 0.40623
$ws.Range("I3").Value = 0.1064658135528677
$ws.Range("J3").Value = 0.1064658135528677
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.220039333333333
$ws.Range("N3").Value = 12.660118
$ws.Range("O3").Value = 0.4495063021722917
$ws.Range("P3").Value = 0.4495063021722917
$ws.Range("Q3").Value = 0.5714355261266667
$ws.Range("R3").Value = 5.14291973514
$ws.Range("S3").Value = 0.04785705415791421
$ws.Range("T3").Value = 0.04785705415791421

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.13541
$ws.Range("H4").Value = 0.40623
$ws.Range("I4").Value = 0.1064658135528677
$ws.Range("J4").Value = 0.1064658135528677
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01908533333333333
$ws.Range("N4").Value = 0.057256
$ws.Range("O4").Value = 0.002032914135332446
$ws.Range("P4").Value = 0.002032914135332446
$ws.Range("Q4").Value = 0.002584344986666666
$ws.Range("R4").Value = 0.02325910488
$ws.Range("S4").Value = 0.0002164358573012934
$ws.Range("T4").Value = 0.0002164358573012934

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.13541
$ws.Range("H5").Value = 0.40623
$ws.Range("I5").Value = 0.1064658135528677
$ws.Range("J5").Value = 0.1064658135528677
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.100466333333333
$ws.Range("N5").Value = 15.301399
$ws.Range("O5").Value = 0.5432868226467401
$ws.Range("P5").Value = 0.5432868226467402
$ws.Range("Q5").Value = 0.6906541461966667
$ws.Range("R5").Value = 6.21588731577
$ws.Range("S5").Value = 0.05784147356563773
$ws.Range("T5").Value = 0.05784147356563774

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gdnf"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.136453666666667
$ws.Range("H6").Value = 3.409361
$ws.Range("I6").Value = 0.8935341864471323
$ws.Range("J6").Value = 0.8935341864471323
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.048574
$ws.Range("N6").Value = 0.145722
$ws.Range("O6").Value = 0.005173961045635648
$ws.Range("P6").Value = 0.005173961045635648
$ws.Range("Q6").Value = 0.05520210040466667
$ws.Range("R6").Value = 0.496818903642
$ws.Range("S6").Value = 0.004623111073621203
$ws.Range("T6").Value = 0.004623111073621203

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gdnf"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.136453666666667
$ws.Range("H7").Value = 3.409361
$ws.Range("I7").Value = 0.8935341864471323
$ws.Range("J7").Value = 0.8935341864471323
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.220039333333333
$ws.Range("N7").Value = 12.660118
$ws.Range("O7").Value = 0.4495063021722917
$ws.Range("P7").Value = 0.4495063021722917
$ws.Range("Q7").Value = 4.795879173844222
$ws.Range("R7").Value = 43.162912564598
$ws.Range("S7").Value = 0.4016492480143775
$ws.Range("T7").Value = 0.4016492480143775

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gdnf"
$ws.Range("C8").Value = "Gfra1"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.136453666666667
$ws.Range("H8").Value = 3.409361
$ws.Range("I8").Value = 0.8935341864471323
$ws.Range("J8").Value = 0.8935341864471323
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01908533333333333
$ws.Range("N8").Value = 0.057256
$ws.Range("O8").Value = 0.002032914135332446
$ws.Range("P8").Value = 0.002032914135332446
$ws.Range("Q8").Value = 0.02168959704622222
$ws.Range("R8").Value = 0.195206373416
$ws.Range("S8").Value = 0.001816478278031153
$ws.Range("T8").Value = 0.001816478278031153

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gdnf"
$ws.Range("C9").Value = "Gfra1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.136453666666667
$ws.Range("H9").Value = 3.409361
$ws.Range("I9").Value = 0.8935341864471323
$ws.Range("J9").Value = 0.8935341864471323
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.100466333333333
$ws.Range("N9").Value = 15.301399
$ws.Range("O9").Value = 0.5432868226467401
$ws.Range("P9").Value = 0.5432868226467402
$ws.Range("Q9").Value = 5.796443666226556
$ws.Range("R9").Value = 52.167992996039
$ws.Range("S9").Value = 0.4854453490811024
$ws.Range("T9").Value = 0.4854453490811025

